$wb = $excel.ActiveWorkbook

# The workbook tracks the same set of events in both the "展览" sheet and
# the combined "全部类型" sheet. Update the "想去人数" (F) and "最低票价" (G)
# figures in both places so they stay in sync.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 1148
    $ws.Range("F6").Value = 141
    $ws.Range("F10").Value = 5200
    $ws.Range("G10").Value = 70
    $ws.Range("F11").Value = 4780
}
